$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (G1 = "0", H1 = "1") as text, matching the
# style already used by the other header cells (bold, bordered, centered).
$ws.Range("G1").Formula = "'0"
$ws.Range("H1").Formula = "'1"

$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data values for the two new columns
$ws.Range("G2").Value = 0
$ws.Range("H3").Value = 1
